$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 55 data ---------------------------------------------------
# Columns (per header row 3): A=Tester B=Algorithm C=Optimization D=#Features
# E=Threshold F=Imputer G=Standardized H=PCA I=Xval Folds J=Local AUC
# K=Kaggle AUC L=Position M=Notes

$ws.Range("A55").Value = "Nanda"
$ws.Range("B55").Value = "GBM, RF"
$ws.Range("C53").Copy($ws.Range("C55"))
$ws.Range("C55").Value = "GBM - BayesOpt + Manual,`nRF - BayesOpt"
$ws.Range("D55").Value = 10
$ws.Range("E55").Value = 0.5

# F55 : "none" (plain text, matches existing style of column F)
$ws.Range("F55").Value = "none"

# G55 / H55 : literal text "False" (not a boolean) - copy from a cell that
# already stores "False" as a shared string so the type is preserved.
$ws.Range("G53").Copy($ws.Range("G55"))
$ws.Range("H53").Copy($ws.Range("H55"))

# I55 : numeric value 5, keep the numeric style used by the rest of the column
$ws.Range("I53").Copy($ws.Range("I55"))
$ws.Range("I55").Value = 5

$ws.Range("L55").Value = "#8"
$ws.Range("J55").Value = "0.866418"
$ws.Range("K55").Value = 0.86899800000000005
$ws.Range("M55").Value = "GBM: n_estimators=197, max_depth=5,  min_samples_split=319,  min_samples_leaf=89,  max_features=0.2, random_state=seed`nRF: n_estimators=161, criterion='gini', min_samples_split=223, min_samples_leaf=9, max_features=1, max_depth=14, random_state=seed`n"

# Row height - content now wraps across more lines
$ws.Rows.Item(55).RowHeight = 63.75

# --- View / layout changes ---------------------------------------------
# Widen column K (index 11) from 18 to ~24.14 characters
$ws.Columns.Item(11).ColumnWidth = 23.25

# Update the active selection to match the new state of the sheet
$ws.Range("I51").Select()
